$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 805
$ws1.Range("F6").Value = 656
$ws1.Range("F9").Value = 809
$ws1.Range("F10").Value = 691
$ws1.Range("F15").Value = 936
$ws1.Range("F16").Value = 9886
$ws1.Range("F17").Value = 617
$ws1.Range("F24").Value = 27
$ws1.Range("F26").Value = 484
$ws1.Range("F27").Value = 180
$ws1.Range("F28").Value = 105
$ws1.Range("F29").Value = 273
$ws1.Range("F32").Value = 67

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 129
$ws2.Range("F21").Value = 1
$ws2.Range("F23").Value = 11

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 805
$ws4.Range("F10").Value = 656
$ws4.Range("F13").Value = 129
$ws4.Range("F14").Value = 809
$ws4.Range("F15").Value = 691
$ws4.Range("F19").Value = 936
$ws4.Range("F20").Value = 9886
$ws4.Range("F22").Value = 617
$ws4.Range("F27").Value = 27
$ws4.Range("F28").Value = 484
$ws4.Range("F29").Value = 180
$ws4.Range("F33").Value = 105
$ws4.Range("F36").Value = 273
$ws4.Range("F39").Value = 67
$ws4.Range("F49").Value = 11

$wb.Save()
